# Issue #13: Permitir que en los ficheros de metadatos dos columnas se puedan
# relacionar para crear SKOS jerarquicos.
#
# The metadata sheet has a header row (row 1, human-readable labels) followed
# by several "property" rows (row 2 onward) that describe each column.
# This change inserts a new property row right after the header that carries
# the raw/API field name for each column (lower-case, hyphenated), pushing
# the previously-existing property rows down by one. The stray trailing
# "mapping-ano.xlsx" value (which lived alone in column E a few rows below
# the real data) is no longer part of the table and is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts the old rows 2-5 down to 3-6).
$ws.Rows.Item(2).Insert()

# Populate the new row with the field/API names for each column.
$ws.Range("A2").Value = "ccaa-nombre"
$ws.Range("B2").Value = "siglas"
$ws.Range("C2").Value = "ccaa-codigo"
$ws.Range("D2").Value = "diputados"
$ws.Range("E2").Value = "ano"
$ws.Range("F2").Value = "votos"

# The old row 5 (now shifted to row 6) only ever held a stray leftover value
# ("mapping-ano.xlsx" in column E) that isn't part of the table any more -
# remove that row entirely so the sheet ends with row 5.
$ws.Rows.Item(6).Delete()
